# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-level holdings detail) right
#    before the "总计" (summary) sheet, mirroring the layout used by the
#    other quarterly sheets (2020-Q4 .. 2021-Q4).
# 2. Insert a new summary row at the top of the "总计" sheet's data for
#    2022-Q1 (35 funds held, 35.73 billion yuan total market value).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet, positioned right before "总计" ----
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Borrow the header-row (B1:H1) and column-A formatting from an existing
# quarterly sheet so the new sheet matches the established look.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A36").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundData = @(
  @(0, "166019", "中欧价值智选回报混合A", "156.17", "94.14", "7.48", "11.6815", 2),
  @(1, "011056", "博时汇兴回报一年持有期灵活配置混合", "107.57", "67.02", "4.40", "4.7331", 5),
  @(2, "013220", "中欧新兴价值一年持有混合A", "63.33", "94.47", "6.37", "4.0341", 4),
  @(3, "004235", "中欧价值智选回报混合C", "36.40", "94.14", "7.48", "2.7227", 2),
  @(4, "001887", "中欧价值智选回报混合E", "20.77", "94.14", "7.48", "1.5536", 2),
  @(5, "000991", "工银瑞信战略转型主题股票A", "48.06", "88.32", "2.44", "1.1727", 10),
  @(6, "011177", "博时汇融回报一年持有期混合A", "39.26", "65.35", "2.85", "1.1189", 6),
  @(7, "004848", "中欧睿泓定期开放灵活配置混合", "23.30", "59.08", "4.80", "1.1184", 3),
  @(8, "013221", "中欧新兴价值一年持有混合C", "16.75", "94.47", "6.37", "1.0670", 4),
  @(9, "519690", "交银稳健配置混合A", "16.44", "88.28", "6.47", "1.0637", 2),
  @(10, "960017", "交银稳健配置混合H", "16.44", "88.28", "6.47", "1.0637", 2),
  @(11, "003378", "泰康策略优选灵活配置混合", "19.86", "81.29", "3.13", "0.6216", 6),
  @(12, "000596", "前海开源中证军工指数A", "14.48", "93.38", "4.27", "0.6183", 8),
  @(13, "012568", "天弘高端制造混合型证券投资基金A", "5.95", "91.15", "8.13", "0.4837", 3),
  @(14, "005977", "中信保诚至兴灵活配置混合A", "7.64", "89.13", "5.35", "0.4087", 4),
  @(15, "010874", "泰康品质生活混合A", "13.17", "81.43", "3.00", "0.3951", 7),
  @(16, "001479", "中邮风格轮动灵活配置混合", "9.45", "62.17", "3.50", "0.3308", 5),
  @(17, "168501", "北信瑞丰产业升级多策略混合", "4.42", "94.11", "5.87", "0.2595", 1),
  @(18, "002199", "前海开源中证军工指数C", "5.45", "93.38", "4.27", "0.2327", 8),
  @(19, "550009", "信诚中小盘混合", "5.23", "87.68", "3.65", "0.1909", 8),
  @(20, "005014", "泰康景泰回报混合A", "11.64", "27.21", "1.34", "0.1560", 8),
  @(21, "011473", "工银瑞信战略转型主题股票C", "5.61", "88.32", "2.44", "0.1369", 10),
  @(22, "010875", "泰康品质生活混合C", "4.39", "81.43", "3.00", "0.1317", 7),
  @(23, "005978", "中信保诚至兴灵活配置混合C", "1.83", "89.13", "5.35", "0.0979", 4),
  @(24, "011927", "博时汇誉回报灵活配置混合型证券投资基金A", "1.30", "68.12", "5.46", "0.0710", 3),
  @(25, "012569", "天弘高端制造混合型证券投资基金C", "0.82", "91.15", "8.13", "0.0667", 3),
  @(26, "009364", "工银瑞信科技创新6个月定期开放混合A", "2.14", "69.85", "2.91", "0.0623", 10),
  @(27, "000535", "长盛航天海工装备灵活配置混合", "2.27", "60.06", "2.33", "0.0529", 10),
  @(28, "164826", "工银瑞信创业板两年定期开放混合A", "2.02", "72.79", "2.57", "0.0519", 10),
  @(29, "009365", "工银瑞信科技创新6个月定期开放混合C", "0.36", "69.85", "2.91", "0.0105", 10),
  @(30, "005015", "泰康景泰回报混合C", "0.63", "27.21", "1.34", "0.0084", 8),
  @(31, "011928", "博时汇誉回报灵活配置混合型证券投资基金C", "0.13", "68.12", "5.46", "0.0071", 3),
  @(32, "011178", "博时汇融回报一年持有期混合C", "0.16", "65.35", "2.85", "0.0046", 6),
  @(33, "010889", "工银瑞信创业板两年定期开放混合C", "0.16", "72.79", "2.57", "0.0041", 10),
  @(34, "002194", "北信瑞丰稳定增强偏债混合", "0.04", "22.00", "1.33", "0.0005", 6)
)
for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# --- 2. Insert a new 2022-Q1 row at the top of the "总计" summary table --
# Re-resolve the "总计" sheet by name: the handle grabbed before
# Worksheets.Add() now points at the freshly-inserted sheet instead (sheet
# handles here track a position, and the new sheet took over slot 6).
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 35
$totalSheet.Range("D2").Value = 35.73
